# Remove the "NABS x.xx - " prefix from every episode title in column A.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$titles = @(
    "The Revolution",
    "Good Night, Sweet Ladies",
    "Random Ghosts",
    "The Lights of Skaro",
    "The Pyramid of Sutekh",
    "The Vaults of Osiris",
    "The Eye of Horus",
    "The Tears of Isis",
    "The Library in the Body",
    "Planet X",
    "The Very Dark Thing",
    "The Emporium at the End",
    "The City and the Clock",
    "Asking For a Friend",
    "Truant",
    "The True Saviour of the Universe",
    "Pride of the Lampian",
    "Clear History",
    "Dead and Breakfast",
    "Burrowed Time",
    "Have I Told You Lately",
    "The Undying Truth",
    "Inertia",
    "Gallifrey"
)

for ($i = 0; $i -lt $titles.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $titles[$i]
}

# Widen column A to fit the (still fairly long) titles.
$ws.Columns.Item(1).ColumnWidth = 27

# Mirror the end-of-edit navigation/selection state: the user scrolled down
# the sheet and selected everything (Ctrl+A) before saving.
$ws.Range("C192").Select() | Out-Null
$ws.Cells.Select() | Out-Null
